# "Add all svm data": turn the single "Train" confusion-matrix sheet into
# two sheets - a brand-new "Train" sheet (the per-fold data) in front, and
# the existing data renamed to "Validation" right after it.

$wb = $excel.ActiveWorkbook

# Duplicate the current (only) sheet so the copy keeps identical styles /
# page setup; Excel places the copy immediately before the source, i.e. at
# index 1, pushing the original to index 2.
$wb.Worksheets.Item(1).Copy($wb.Worksheets.Item(1))

$trainWs = $wb.Worksheets.Item(1)
$validationWs = $wb.Worksheets.Item(2)

$validationWs.Name = "Validation"
$trainWs.Name = "Train"

# --- Validation sheet: just move the selection, content is untouched ---
$validationWs.Select()
$validationWs.Range("B13:L14").Select()

# --- Train sheet: replace the copied values with the new dataset ---
$ws = $trainWs

$ws.Range("C2").Value = 0
$ws.Range("D2").Value = "affirmative"
$ws.Range("E2").Value = "conditional"
$ws.Range("F2").Value = "doubt_question"
$ws.Range("G2").Value = "emphasis"
$ws.Range("H2").Value = "negative"
$ws.Range("I2").Value = "relative"
$ws.Range("J2").Value = "topics"
$ws.Range("K2").Value = "wh_question"
$ws.Range("L2").Value = "yn_question"

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 2532
$ws.Range("D3").Value = 85
$ws.Range("E3").Value = 63
$ws.Range("F3").Value = 56
$ws.Range("G3").Value = 39
$ws.Range("H3").Value = 94
$ws.Range("I3").Value = 69
$ws.Range("J3").Value = 44
$ws.Range("K3").Value = 122
$ws.Range("L3").Value = 75

$ws.Range("B4").Value = "affirmative"
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 46
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0

$ws.Range("B5").Value = "conditional"
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 100
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0

$ws.Range("B6").Value = "doubt_question"
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 120
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 1

$ws.Range("B7").Value = "emphasis"
$ws.Range("C7").Value = 4
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 75
$ws.Range("H7").Value = 2
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0

$ws.Range("B8").Value = "negative"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 75
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0

$ws.Range("B9").Value = "relative"
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 99
$ws.Range("J9").Value = 4
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 1

$ws.Range("B10").Value = "topics"
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 70
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0

$ws.Range("B11").Value = "wh_question"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 42
$ws.Range("L11").Value = 0

$ws.Range("B12").Value = "yn_question"
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 2
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 2
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 109

# Row 13 ("hit") / row 14 ("miss") keep the same formula shapes as the
# original sheet, just re-entered so they recalc against the new grid.
$ws.Range("B13").Value = "hit"
$ws.Range("C13").Formula = "=C3"
$ws.Range("D13").Formula = "=D4"
$ws.Range("E13").Formula = "=E5"
$ws.Range("F13").Formula = "=F6"
$ws.Range("G13").Formula = "=G7"
$ws.Range("H13").Formula = "=H8"
$ws.Range("I13").Formula = "=I9"
$ws.Range("J13").Formula = "=J10"
$ws.Range("K13").Formula = "=K11"
$ws.Range("L13").Formula = "=L12"

$ws.Range("B14").Value = "miss"
$ws.Range("C14").Formula = "=SUM(C3:C12) - C13"
$ws.Range("D14:L14").Formula = "=SUM(D3:D12) - D13"

# The source sheet had some blank formatted cells trailing the totals row
# (O13:X13 / O14:X14); recreate them purely as formatting (no value) by
# copying the "hit"/"miss" row style across, matching the widened B2:X14
# used range.
$ws.Range("C13").Copy()
$ws.Range("O13:X14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Drop the leftover rows (15-26) that belonged to the old, bigger sheet -
# the new Train sheet only spans down to row 14.
$ws.Range("B15:N26").Clear()

# Clear the inherited "B column" emphasis formatting - the new sheet's
# label column (and the header row) carries no special style.
$ws.Range("B2").ClearContents()
$ws.Range("B2").ClearFormats()
$ws.Range("B3:B12").ClearFormats()

$ws.Select()
$ws.Range("R19").Select()

